$wb = $excel.ActiveWorkbook

# Rename "Sheet1" to "Nädal 5"
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Nädal 5"

# Row 9: add Stop time and Delta Time, change Activity/Comments
$ws.Range("D9").Value = 0.9375
$ws.Range("F9").Value = 90
$ws.Range("G9").Value = "vead"
$ws.Range("H9").Value = "vigade parandus"

# Row 10: fill in Date, Start, Stop, Delta Time, Activity, Comments
$ws.Range("B10").Value = 43893
$ws.Range("C10").Value = 0.33333333333333331
$ws.Range("D10").Value = 0.39583333333333331
$ws.Range("F10").Value = 90
$ws.Range("G10").Value = "loeng"
$ws.Range("H10").Value = "W6 loeng"

# Row 11: fill in Date, Start, Stop, Delta Time, Activity, Comments, C
$ws.Range("B11").Value = 43893
$ws.Range("C11").Value = 0.77777777777777779
$ws.Range("D11").Value = 0.8125
$ws.Range("F11").Value = 50
$ws.Range("G11").Value = "video"
$ws.Range("H11").Value = "HW5, video 7 lõpuni"
$ws.Range("I11").Value = "x"

# Update selection to H12
$ws.Range("H12").Select()
